# Weekly Time Tracker - add "Work Day" / name header, update hours for two days

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells: D3 "Work Day" label, D1 the name, D2 just a centered spacer cell.
# Set "Work Day" (D3) before "Scott McNulty" (D1) so the shared-string table grows
# in that order (index 18 = "Work Day", index 19 = "Scott McNulty").
$ws.Range("D3").Value = "Work Day"
$ws.Range("D1").Value = "Scott McNulty"

$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4108
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108

$ws.Range("D2").HorizontalAlignment = -4108

# Fill in missing hours for the "Week 2" block (row 11)
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 2

# Fill in hours for the "Week 4" block (row 16)
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 7

# Update the active selection to match the author's last-edited cell
$ws.Range("M14").Select() | Out-Null
